$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "New category"
$ws.Range("C1").Value = "Taco category"
$ws.Range("B2").Value = "Non-recyclable"
$ws.Range("C2").Value = "Aluminium foil"
$ws.Range("B3").Value = "E-Waste"
$ws.Range("C3").Value = "Battery"
$ws.Range("C4").Value = "Aluminium blister pack"
$ws.Range("B5").Value = "Non-recyclable"
$ws.Range("C5").Value = "Carded blister pack"
$ws.Range("B6").Value = "Plastic"
$ws.Range("C6").Value = "Other plastic bottle"
$ws.Range("B7").Value = "Plastic"
$ws.Range("C7").Value = "Clear plastic bottle"
$ws.Range("B8").Value = "Glass"
$ws.Range("C8").Value = "Glass bottle"
$ws.Range("C9").Value = "Plastic bottle cap"
$ws.Range("B10").Value = "Metal"
$ws.Range("C10").Value = "Metal bottle cap"
$ws.Range("B11").Value = "Glass"
$ws.Range("C11").Value = "Broken glass"
$ws.Range("B12").Value = "Non-recyclable"
$ws.Range("C12").Value = "Food Can"
$ws.Range("B13").Value = "Metal"
$ws.Range("C13").Value = "Aerosol"
$ws.Range("C15").Value = "Toilet tube"
$ws.Range("C16").Value = "Other carton"
$ws.Range("B17").Value = "Paper"
$ws.Range("C17").Value = "Egg carton"
$ws.Range("B18").Value = "Paper"
$ws.Range("C18").Value = "Drink carton"
$ws.Range("B19").Value = "Paper"
$ws.Range("C19").Value = "Corrugated carton"
$ws.Range("B20").Value = "Paper"
$ws.Range("C20").Value = "Meal carton"
$ws.Range("B21").Value = "Paper"
$ws.Range("C21").Value = "Pizza box"
$ws.Range("B22").Value = "Paper"
$ws.Range("C22").Value = "Paper cup"
$ws.Range("B23").Value = "Plastic"
$ws.Range("C23").Value = "Disposable plastic cup"
$ws.Range("B24").Value = "Non-recyclable"
$ws.Range("C24").Value = "Foam cup"
$ws.Range("B25").Value = "Non-recyclable"
$ws.Range("C25").Value = "Glass cup"
$ws.Range("B26").Value = "Plastic"
$ws.Range("C26").Value = "Other plastic cup"
$ws.Range("B27").Value = "Organic"
$ws.Range("C27").Value = "Food waste"
$ws.Range("B28").Value = "Glass"
$ws.Range("C28").Value = "Glass jar"
$ws.Range("B29").Value = "Plastic"
$ws.Range("C29").Value = "Plastic lid"
$ws.Range("B30").Value = "Metal"
$ws.Range("C30").Value = "Metal lid"
$ws.Range("B32").Value = "Paper"
$ws.Range("C32").Value = "Magazine paper"
$ws.Range("B33").Value = "Non-recyclable"
$ws.Range("C33").Value = "Tissues"
$ws.Range("B34").Value = "Paper"
$ws.Range("C34").Value = "Wrapping paper"
$ws.Range("B35").Value = "Paper"
$ws.Range("C35").Value = "Normal paper"
$ws.Range("B36").Value = "Paper"
$ws.Range("B37").Value = "Non-recyclable"
$ws.Range("C37").Value = "Plastified paper bag"
$ws.Range("B38").Value = "Plastic"
$ws.Range("C38").Value = "Plastic film"
$ws.Range("B39").Value = "Plastic"
$ws.Range("C39").Value = "Six pack rings"
$ws.Range("B40").Value = "Non-recyclable"
$ws.Range("C40").Value = "Garbage bag"
$ws.Range("C41").Value = "Other plastic wrapper"
$ws.Range("C42").Value = "Single-use carrier bag"
$ws.Range("B43").Value = "Non-recyclable"
$ws.Range("C43").Value = "Polypropylene bag"
$ws.Range("C44").Value = "Crisp packet"
$ws.Range("C45").Value = "Spread tub"
$ws.Range("B46").Value = "Plastic"
$ws.Range("C46").Value = "Tupperware"
$ws.Range("B47").Value = "Plastic"
$ws.Range("C47").Value = "Disposable food container"
$ws.Range("B48").Value = "Non-recyclable"
$ws.Range("C48").Value = "Foam food container"
$ws.Range("B49").Value = "Plastic"
$ws.Range("C49").Value = "Other plastic container"
$ws.Range("B50").Value = "Plastic"
$ws.Range("C50").Value = "Plastic glooves"
$ws.Range("B51").Value = "Plastic"
$ws.Range("C51").Value = "Plastic utensils"
$ws.Range("B52").Value = "Metal"
$ws.Range("C52").Value = "Pop tab"
$ws.Range("B53").Value = "Non-recyclable"
$ws.Range("C53").Value = "Rope & strings"
$ws.Range("B54").Value = "Metal"
$ws.Range("C54").Value = "Scrap metal"
$ws.Range("C55").Value = "Shoe"
$ws.Range("C56").Value = "Squeezable tube"
$ws.Range("B57").Value = "Plastic"
$ws.Range("C57").Value = "Plastic straw"
$ws.Range("C58").Value = "Paper straw"
$ws.Range("B59").Value = "Non-recyclable"
$ws.Range("C59").Value = "Styrofoam piece"
$ws.Range("B61").Value = "Non-recyclable"
$ws.Range("C61").Value = "Cigarette"
